$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Machines Current Steps/mm" (B) and "Actual Distance Moved (mm)" (D)
# values for the E-step rows (E0-E3, rows 7-10). The commanded distance (C)
# stays 100 for all of them, and E/F are formula-driven so they recalc
# automatically.

$ws.Range("B7").Value = 1070
$ws.Range("D7").Value = 96

$ws.Range("B8").Value = 1070
$ws.Range("D8").Value = 100.5

$ws.Range("B9").Value = 1070
$ws.Range("D9").Value = 102

$ws.Range("B10").Value = 1070
$ws.Range("D10").Value = 95

# Update selection / view to match the saved state in the target workbook.
$ws.Range("E30").Select()
